$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the "Unavailability" note in B2 (was "1->4")
$ws.Range("B2").ClearContents()

# Update "Parallel exams" (column C) values to reflect added dependencies
# for different proms but same professor.
$ws.Range("C2").Value = 2
$ws.Range("C3").Value = 2
$ws.Range("C4").Value = 2
$ws.Range("C6").Value = 1
$ws.Range("C8").Value = 2
$ws.Range("C9").Value = 1
$ws.Range("C10").Value = 2

# Move the active selection from C10 to C9
$ws.Range("C9").Select()
